$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows (3,4,5,7,8,9,10,11) have been rearranged: each row now carries
# the values that previously belonged to a different row (row 6 is untouched).
# Columns touched: D (Fecha), L (Calidad), M (Volumen), N (Precio minimo),
# O (Precio maximo), P (Precio promedio ponderado), Q (Unidad de
# comercializacion), R (Origen), S (Precio $/Kg), T (Kg / unidad)

$rows = @{
    3  = @{ D = 44511; L = "Primera"; M = 45; N = 28000; O = 28000; P = 28000; Q = "`$/bandeja 10 kilos"; R = "Provincia de Los Andes"; S = 2800; T = 10 }
    4  = @{ D = 44511; L = "Primera"; M = 45; N = 3200;  O = 3200;  P = 3200;  Q = "`$/bandeja 10 kilos"; R = "Provincia de Quillota";  S = 320;  T = 10 }
    5  = @{ D = 44166; L = "Segunda"; M = 20; N = 12000; O = 12000; P = 12000; Q = "`$/caja 18 kilos";    R = "La Ligua";                S = 667;  T = 18 }
    7  = @{ D = 44483; L = "Primera"; M = 35; N = 10000; O = 10000; P = 10000; Q = "`$/bandeja 5 kilos";  R = "Provincia de Quillota";  S = 2000; T = 5  }
    8  = @{ D = 44503; L = "Primera"; M = 50; N = 28000; O = 28000; P = 28000; Q = "`$/bandeja 10 kilos"; R = "Provincia de Quillota";  S = 2800; T = 10 }
    9  = @{ D = 44519; L = "Primera"; M = 30; N = 28000; O = 28000; P = 28000; Q = "`$/bandeja 10 kilos"; R = "Provincia de Quillota";  S = 2800; T = 10 }
    10 = @{ D = 44496; L = "Primera"; M = 55; N = 28000; O = 28000; P = 28000; Q = "`$/bandeja 10 kilos"; R = "Provincia de Quillota";  S = 2800; T = 10 }
    11 = @{ D = 44515; L = "Primera"; M = 80; N = 28000; O = 28000; P = 28000; Q = "`$/bandeja 10 kilos"; R = "Provincia de Los Andes"; S = 2800; T = 10 }
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Cells.Item($r, 4).Value = $vals.D    # D - Fecha
    $ws.Cells.Item($r, 12).Value = $vals.L   # L - Calidad
    $ws.Cells.Item($r, 13).Value = $vals.M   # M - Volumen
    $ws.Cells.Item($r, 14).Value = $vals.N   # N - Precio minimo
    $ws.Cells.Item($r, 15).Value = $vals.O   # O - Precio maximo
    $ws.Cells.Item($r, 16).Value = $vals.P   # P - Precio promedio ponderado
    $ws.Cells.Item($r, 17).Value = $vals.Q   # Q - Unidad de comercializacion
    $ws.Cells.Item($r, 18).Value = $vals.R   # R - Origen
    $ws.Cells.Item($r, 19).Value = $vals.S   # S - Precio $/Kg
    $ws.Cells.Item($r, 20).Value = $vals.T   # T - Kg / unidad
}
